$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking need NumberFormat forced to Text
# first, otherwise Excel auto-converts the literal text into a Number and the
# original formatting (e.g. trailing zeros, thousands-dot grouping) is lost.

$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D15", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptos data (price + 1h volume refresh, a couple of
# coin rows also got re-ordered / relabeled by the upstream feed).

$ws.Range("D2").Value = '26.084.78'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.646.21'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '215.19'
$ws.Range("E5").Value = '  +2.68%  '
$ws.Range("D6").Value = '0.5224'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.2605'
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").Value = '0.06314'
$ws.Range("E9").Value = '  -0.48%  '
$ws.Range("D10").Value = '20.77'
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("D11").Value = '0.07684'
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("D12").Value = '1.647.19'
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("D13").Value = '4.420'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").Value = '1.868.66'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").Value = '0.5569'
$ws.Range("E15").Value = '  +1.24%  '
$ws.Range("D16").Value = '0.0₅8194'
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("D17").Value = '65.21'
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = '26.101.25'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("D20").Value = '4.732'
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").Value = '188.87'
$ws.Range("E21").Value = '  +1.18%  '
$ws.Range("D22").Value = '10.22'
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D23").Value = '6.202'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").Value = '145.91'
$ws.Range("E25").Value = '  -2.41%  '
$ws.Range("D26").Value = '7.437'
$ws.Range("E26").Value = '  -0.87%  '
$ws.Range("D27").Value = '0.1208'
$ws.Range("E27").Value = '  -3.26%  '
$ws.Range("D28").Value = '15.87'
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").Value = '1.392'
$ws.Range("E29").Value = '  +2.83%  '
$ws.Range("D30").Value = '0.05922'
$ws.Range("D31").Value = '1.263'
$ws.Range("E31").Value = '  -0.96%  '
$ws.Range("D32").Value = '3.442'
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("D33").Value = '3.407'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").Value = '1.654'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").Value = '0.9849'
$ws.Range("E35").Value = '  -2.10%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.394'
$ws.Range("E36").Value = '  -0.51%  '
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").Value = '2.759'
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").Value = '0.5676'
$ws.Range("E38").Value = '  -5.81%  '
$ws.Range("D39").Value = '0.01617'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").Value = '5.777'
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '1.028.19'
$ws.Range("E43").Value = '  -7.63%  '
$ws.Range("D44").Value = '100.24'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '1.795.69'
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₈109'
$ws.Range("E46").Value = '  -1.73%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '55.99'
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").Value = '0.9997'
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("D49").Value = '8.071'
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").Value = '0.05187'
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("D51").Value = '0.4221'
$ws.Range("E51").Value = '  -0.56%  '
